# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# and wrap the sheet's data range in a table ("Table1"), then freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the 21 header cells in row 1 (A1:U1).
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range (A1:U64) into a table so the renamed headers
#    become the table's column names.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split below row 1, top-left cell of the
#    scrolling pane is A2) - select the cell below the freeze line
#    first, matching the native Excel "Freeze Panes" gesture.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
